$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '69.349.97'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.685.29'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '678.60'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.15'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.07%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('E10').Value = '  -3.89%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.437'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('E12').Value = '  -3.11%  '
$ws.Range('D13').Value = '4.308.75'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('E14').Value = '  -3.42%  '
$ws.Range('D15').Value = '3.690.10'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '69.289.90'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '16.02'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '468.47'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('D24').Value = '3.831.91'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  -5.37%  '
$ws.Range('E27').Value = '  -4.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.11'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.25%  '
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('E30').Value = '  -3.35%  '
$ws.Range('E31').Value = '  -3.39%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.94'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.71%  '
$ws.Range('E34').Value = '  -4.41%  '
$ws.Range('D35').Value = '3.675.37'
$ws.Range('E35').Value = '  +0.69%  '
$ws.Range('E36').Value = '  -4.59%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '8.28'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.22'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.39%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('E40').Value = '  -3.84%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0905'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.96%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '170.18'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.942'
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '47.69'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.71'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.38%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.000279'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.93%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '27.98'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -6.30%  '
$ws.Range('E50').Value = '  -5.45%  '
$ws.Range('E51').Value = '  -2.67%  '
